# Generate Report for Handback
# Refresh the handoff/handback timestamps for the "346df058-44a9-4cd4-a2e4-958e3d85f3eb"
# entry across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" for the 346df058... row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-06 06:52:53"

# --- zh-cn sheet: Correspond Handoff / Handback DateTime for the 346df058... row ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-06 06:52:47"
$wsZhCn.Range("K3").Value = "2016-09-06 06:53:13"

# --- de-de sheet: Correspond Handoff / Handback DateTime for the 346df058... row ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-06 06:52:53"
$wsDeDe.Range("K3").Value = "2016-09-06 06:53:21"
